# Updated symbol list on Sat Jan 14 15:16:43 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures for the
# coin rows on Sheet1, matching the latest scrape of coinranking.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell address paired with its new textual value. The
# values are written with a leading apostrophe so Excel stores them as
# text (preserving formatting such as trailing zeros, thousand
# separators and the trailing "%" sign) instead of coercing them into
# numbers, matching how the sheet already represents these columns.
$updates = @(
    @{ Cell = "D2"; Value = "303.50" },
    @{ Cell = "E2"; Value = "5.90%" },
    @{ Cell = "D3"; Value = "31.96" },
    @{ Cell = "E3"; Value = "8.72%" },
    @{ Cell = "D4"; Value = "5.305" },
    @{ Cell = "E4"; Value = "4.54%" },
    @{ Cell = "D5"; Value = "0.07458" },
    @{ Cell = "E5"; Value = "10.11%" },
    @{ Cell = "D6"; Value = "7.835" },
    @{ Cell = "E6"; Value = "6.36%" },
    @{ Cell = "D7"; Value = "3.807" },
    @{ Cell = "E7"; Value = "10.65%" },
    @{ Cell = "D8"; Value = "1.452" },
    @{ Cell = "E8"; Value = "5.62%" },
    @{ Cell = "D9"; Value = "0.9195" },
    @{ Cell = "E9"; Value = "2.16%" },
    @{ Cell = "D10"; Value = "0.01743" },
    @{ Cell = "E10"; Value = "2,606.94%" },
    @{ Cell = "D11"; Value = "0.1692" },
    @{ Cell = "E11"; Value = "6.60%" },
    @{ Cell = "D12"; Value = "0.07709" },
    @{ Cell = "E12"; Value = "11.83%" },
    @{ Cell = "D13"; Value = "0.08016" },
    @{ Cell = "E13"; Value = "5.46%" },
    @{ Cell = "D14"; Value = "0.03028" },
    @{ Cell = "E14"; Value = "3.68%" },
    @{ Cell = "D15"; Value = "0.09881" },
    @{ Cell = "E15"; Value = "9.87%" },
    @{ Cell = "E16"; Value = "-5.09%" },
    @{ Cell = "D17"; Value = "0.04579" },
    @{ Cell = "E17"; Value = "2.21%" },
    @{ Cell = "D18"; Value = "0.006224" },
    @{ Cell = "E18"; Value = "-0.49%" },
    @{ Cell = "E19"; Value = "0.69%" },
    @{ Cell = "D20"; Value = "2.228" },
    @{ Cell = "E20"; Value = "-0.03%" },
    @{ Cell = "E21"; Value = "3.16%" },
    @{ Cell = "D22"; Value = "0.1345" },
    @{ Cell = "E22"; Value = "1.69%" },
    @{ Cell = "D23"; Value = "4.508" },
    @{ Cell = "E23"; Value = "12.73%" },
    @{ Cell = "D24"; Value = "0.1623" },
    @{ Cell = "E24"; Value = "4.28%" },
    @{ Cell = "E25"; Value = "1.29%" },
    @{ Cell = "D26"; Value = "0.004410" },
    @{ Cell = "E26"; Value = "1.04%" },
    @{ Cell = "D27"; Value = "0.0001401" },
    @{ Cell = "E27"; Value = "20.27%" },
    @{ Cell = "D28"; Value = "0.0001742" },
    @{ Cell = "E28"; Value = "-1.46%" },
    @{ Cell = "D40"; Value = "0.04518" },
    @{ Cell = "E40"; Value = "5.74%" },
    @{ Cell = "D41"; Value = "0.007211" },
    @{ Cell = "E41"; Value = "6.36%" },
    @{ Cell = "D42"; Value = "0.1341" },
    @{ Cell = "E42"; Value = "8.14%" },
    @{ Cell = "D43"; Value = "0.002242" },
    @{ Cell = "E43"; Value = "2.34%" },
    @{ Cell = "D44"; Value = "0.01260" },
    @{ Cell = "E44"; Value = "9.50%" },
    @{ Cell = "D45"; Value = "0.00006156" },
    @{ Cell = "E45"; Value = "7.69%" },
    @{ Cell = "D47"; Value = "0.01301" },
    @{ Cell = "E47"; Value = "-0.04%" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
}
